# Apply the itinerary edit: change the 07:30 AM "Arrival / Parking Lot 1"
# entry on Sheet1 to a 09:15 AM "Coffee Time / muh house" entry, and make
# Sheet1 the active/selected sheet (it was Sheet2 before).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update the data row (row 2) on Sheet1: time, title, location.
$ws1.Cells.Item(2, 1).Value = 0.38541666666666669
$ws1.Cells.Item(2, 2).Value = "Coffee Time"
$ws1.Cells.Item(2, 3).Value = "muh house"

# Sheet1 becomes the selected/active sheet; move its selection to B3.
$ws1.Activate()
$ws1.Range("B3").Select()
